$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 551.7911908829242
$ws.Range("D2").Value = 134.8010645558176
$ws.Range("F2").Value = 452
$ws.Range("G2").Value = 508
$ws.Range("H2").Value = 619

$ws.Range("C3").Value = 40.5752808742658
$ws.Range("D3").Value = 4.739821021780309
$ws.Range("G3").Value = 39.89
$ws.Range("H3").Value = 43.15

$ws.Range("C4").Value = 1.41712710083542
$ws.Range("D4").Value = 2.243893152512907
$ws.Range("H4").Value = 1.82

$ws.Range("C5").Value = 324.0963851543303
$ws.Range("D5").Value = 10.18603570960344
$ws.Range("F5").Value = 318.16
$ws.Range("G5").Value = 326
$ws.Range("H5").Value = 332.46

$ws.Range("C6").Value = 20.84400402499239
$ws.Range("D6").Value = 2.285550540142486
$ws.Range("F6").Value = 19.4
$ws.Range("G6").Value = 20.54
$ws.Range("H6").Value = 22.15

$ws.Range("C7").Value = -76.34821987691012
$ws.Range("D7").Value = 22.43616008784592

$ws.Range("C8").Value = 7.75206201747966
$ws.Range("D8").Value = 6.830867396907287

$ws.Range("C9").Value = 9.322889570121452
$ws.Range("D9").Value = 1.688110164882354

$ws.Range("C10").Value = 867.8303416095284
$ws.Range("D10").Value = 0.4610683516698383

$ws.Range("C11").Value = 0.5569473730409753
$ws.Range("D11").Value = 0.5905526589122854

$ws.Range("C12").Value = 22.68879601244939
$ws.Range("D12").Value = 12.27778629887413

$ws.Range("C13").Value = 0.6714604638101701
$ws.Range("D13").Value = 0.7482366910637375

$ws.Range("C14").Value = 1.826063229822386
$ws.Range("D14").Value = 1.666121582407144

$ws.Range("C15").Value = 93.7482198769099
$ws.Range("D15").Value = 22.43616008784592

$ws.Range("C16").Value = -85.65299769868724
$ws.Range("D16").Value = 20.21837856105878
$ws.Range("G16").Value = -85.45410721860875
$ws.Range("H16").Value = -67.95746206410165

$ws.Range("C17").Value = -77.90093568120759
$ws.Range("D17").Value = 24.80926135014304
$ws.Range("F17").Value = -92.29706163635328
$ws.Range("G17").Value = -75.29706163635328
$ws.Range("H17").Value = -57.26572375596102
